$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the text content of row 1 and row 2 (columns A and B), so that
# the "Equipamientos de salud..." pair moves to row 1 (indices 0,1 in the
# shared string table) and "No aplicable" pair moves to row 2 (indices 2,3),
# matching the new shared-strings order in the target file.
$a1 = $ws.Range("A1").Value2
$b1 = $ws.Range("B1").Value2
$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2

$ws.Range("A1").Value = $a2
$ws.Range("B1").Value = $b2
$ws.Range("A2").Value = $a1
$ws.Range("B2").Value = $b1
